# Insert a new data row at row 453 (pushes existing rows 453..482 down to 454..483)
# and populate it with a new "Cilantro" price record, matching the diff:
#   - dimension grows from A1:R482 to A1:R483
#   - old rows 453..482 become 454..483 (shifted down by one), unchanged
#   - new row 453 carries fresh data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 453:482 down to 454:483, leaving a blank row 453 behind.
$ws.Rows.Item(453).Insert()

# Populate the newly inserted row 453 with the new record.
$ws.Cells.Item(453, 1).Value = 3
$ws.Cells.Item(453, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(453, 3).Value = "Coquimbo"
$ws.Cells.Item(453, 4).Value = 44931
$ws.Cells.Item(453, 5).Value = 5
$ws.Cells.Item(453, 6).Value = 100112040
$ws.Cells.Item(453, 7).Value = "Cilantro"
$ws.Cells.Item(453, 8).Value = "Sin especificar"
$ws.Cells.Item(453, 9).Value = "Primera"
$ws.Cells.Item(453, 10).Value = 200
$ws.Cells.Item(453, 11).Value = 4500
$ws.Cells.Item(453, 12).Value = 5000
$ws.Cells.Item(453, 13).Value = 4775
$ws.Cells.Item(453, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(453, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(453, 16).Value = 1592
$ws.Cells.Item(453, 17).Value = 3
$ws.Cells.Item(453, 18).Value = "Hortaliza"
